$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all "VARCHAR(30)" entries (Type column B) to "VARCHAR(50)"
foreach ($row in 2..12) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value2 -eq "VARCHAR(30)") {
        $cell.Value = "VARCHAR(50)"
    }
}

# Restore default view: no frozen/scrolled top-left cell, and select B9:B12
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9:B12").Select()
